$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.356.69'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '1.941.27'
$ws.Range('E3').Value = '  -2.95%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '252.48'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = '0.7185'
$ws.Range('E6').Value = '  -7.39%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.3346'
$ws.Range('E8').Value = '  -4.03%  '
$ws.Range('D9').Value = '28.66'
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('D10').Value = '0.07398'
$ws.Range('E10').Value = '  +5.63%  '
$ws.Range('D11').Value = '0.8167'
$ws.Range('E11').Value = '  -4.13%  '
$ws.Range('D12').Value = '0.08149'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = '1.940.44'
$ws.Range('E13').Value = '  -2.98%  '
$ws.Range('D14').Value = '5.497'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').Value = '95.29'
$ws.Range('E15').Value = '  -5.37%  '
$ws.Range('D16').Value = '14.92'
$ws.Range('E16').Value = '  -3.34%  '
$ws.Range('D17').Value = '30.382.34'
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').Value = '0.000008369'
$ws.Range('E18').Value = '  +5.98%  '
$ws.Range('D19').Value = '254.43'
$ws.Range('E19').Value = '  -7.50%  '
$ws.Range('D20').Value = '5.866'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('D21').Value = '2.196.57'
$ws.Range('E21').Value = '  -2.63%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '6.992'
$ws.Range('E24').Value = '  -1.74%  '
$ws.Range('D25').Value = '9.875'
$ws.Range('E25').Value = '  -1.64%  '
$ws.Range('D26').Value = '161.60'
$ws.Range('E26').Value = '  -1.85%  '
$ws.Range('D27').Value = '2.417'
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('D28').Value = '19.41'
$ws.Range('E28').Value = '  -2.50%  '
$ws.Range('D29').Value = '0.1318'
$ws.Range('E29').Value = '  -10.76%  '
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('D32').Value = '4.489'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('D33').Value = '4.269'
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('D34').Value = '0.05295'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('D35').Value = '1.284'
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('D36').Value = '0.7628'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('D37').Value = '2.749'
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = '0.02000'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').Value = '81.64'
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').Value = '6.592'
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').Value = '0.4576'
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').Value = '2.038'
$ws.Range('E43').Value = '  -5.08%  '
$ws.Range('D44').Value = '0.8472'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = '103.33'
$ws.Range('E46').Value = '  -2.29%  '
$ws.Range('D47').Value = '9.879'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').Value = '7.484'
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('D49').Value = '37.29'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('D50').Value = '0.4203'
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('D51').Value = '0.06056'
$ws.Range('E51').Value = '  +1.31%  '
